# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N ("Late") on the
# "Repayment schedule" sheet, shifting the existing Late / heading /
# Outstanding columns one place to the right, then size the new column
# and make "Repayment schedule" the active/selected sheet & cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column at N - existing N/O/P (Late, heading, Outstanding)
# shift right to O/P/Q.
$ws.Columns("N").Insert()

# Match the authored column width for the newly inserted column.
$ws.Columns("N").ColumnWidth = 10.2

# Make "Repayment schedule" the active sheet/tab with R6 selected.
$ws.Activate()
$ws.Range("R6").Select()
